# Apply the updated project data values and refresh the sheet's selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ProjectDependency (C) / ProjectGroup (D) columns get simplified to 1.
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 1

# Reflect the table range as the active selection on the sheet.
$ws.Range("A1:D5").Select()
